# Apply updated crypto price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B (Coin) and C (Link) swaps for rows 47-48 ---
$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"

# --- Column D (Price) updates ---
# Force text format so numeric-looking strings (e.g. "1.00") are not
# auto-converted to numbers, matching the original inline-string cells.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.685.31"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.411.00"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.49"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.95"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.385"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.993.87"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.410.34"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.762.03"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.13"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.77"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.14"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "387.39"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.32"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.189"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.38"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.99"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.42"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.18"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.94"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "168.57"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.445.61"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.46"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "28.34"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0754"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.785"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.500.56"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.75"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.62"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0263"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.205"

# --- Column E (Volume 1h) updates ---
$ws.Range("E2").Value = "  +1.08%  "
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("E5").Value = "  +1.00%  "
$ws.Range("E6").Value = "  +1.38%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("E10").Value = "  -0.06%  "
$ws.Range("E11").Value = "  -0.66%  "
$ws.Range("E12").Value = "  +0.67%  "
$ws.Range("E13").Value = "  -0.59%  "
$ws.Range("E14").Value = "  +0.47%  "
$ws.Range("E15").Value = "  +0.60%  "
$ws.Range("E16").Value = "  -1.23%  "
$ws.Range("E17").Value = "  +1.07%  "
$ws.Range("E18").Value = "  +0.67%  "
$ws.Range("E19").Value = "  +1.18%  "
$ws.Range("E20").Value = "  +2.59%  "
$ws.Range("E21").Value = "  +0.93%  "
$ws.Range("E22").Value = "  -1.17%  "
$ws.Range("E23").Value = "  -0.47%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("E25").Value = "  -1.45%  "
$ws.Range("E26").Value = "  +3.41%  "
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("E28").Value = "  +1.48%  "
$ws.Range("E29").Value = "  +0.23%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("E31").Value = "  +1.40%  "
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("E33").Value = "  +0.88%  "
$ws.Range("E34").Value = "  +4.20%  "
$ws.Range("E35").Value = "  -0.18%  "
$ws.Range("E36").Value = "  +1.15%  "
$ws.Range("E37").Value = "  +0.76%  "
$ws.Range("E38").Value = "  -0.28%  "
$ws.Range("E39").Value = "  +5.62%  "
$ws.Range("E40").Value = "  -1.72%  "
$ws.Range("E41").Value = "  +0.70%  "
$ws.Range("E42").Value = "  +1.50%  "
$ws.Range("E43").Value = "  +0.42%  "
$ws.Range("E44").Value = "  +3.33%  "
$ws.Range("E45").Value = "  +2.05%  "
$ws.Range("E46").Value = "  -0.86%  "
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("E48").Value = "  -1.38%  "
$ws.Range("E49").Value = "  -0.21%  "
$ws.Range("E50").Value = "  -2.57%  "
$ws.Range("E51").Value = "  -0.70%  "
